$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2
$ws.Range("G2").Value = 113.5893336666667
$ws.Range("H2").Value = 340.768001
$ws.Range("I2").Value = 0.8306211829777892
$ws.Range("J2").Value = 0.8306211829777892
$ws.Range("M2").Value = 14.82137566666667
$ws.Range("N2").Value = 44.464127
$ws.Range("O2").Value = 0.9193095050964014
$ws.Range("P2").Value = 0.9193095050964012
$ws.Range("Q2").Value = 1683.550186000014
$ws.Range("R2").Value = 15151.95167400013
$ws.Range("S2").Value = 0.7635979486458989
$ws.Range("T2").Value = 0.7635979486458987

# Row 3
$ws.Range("G3").Value = 113.5893336666667
$ws.Range("H3").Value = 340.768001
$ws.Range("I3").Value = 0.8306211829777892
$ws.Range("J3").Value = 0.8306211829777892
$ws.Range("O3").Value = 0.0119476607067887
$ws.Range("P3").Value = 0.0119476607067887
$ws.Range("Q3").Value = 21.87999394509678
$ws.Range("R3").Value = 196.919945505871
$ws.Range("S3").Value = 0.009923980070090078
$ws.Range("T3").Value = 0.009923980070090077

# Row 4
$ws.Range("G4").Value = 113.5893336666667
$ws.Range("H4").Value = 340.768001
$ws.Range("I4").Value = 0.8306211829777892
$ws.Range("J4").Value = 0.8306211829777892
$ws.Range("M4").Value = 1.108292
$ws.Range("N4").Value = 3.324876
$ws.Range("O4").Value = 0.06874283419681
$ws.Range("P4").Value = 0.06874283419680997
$ws.Range("Q4").Value = 125.8901497880973
$ws.Range("R4").Value = 1133.011348092876
$ws.Range("S4").Value = 0.05709925426180035
$ws.Range("T4").Value = 0.05709925426180032

# Row 5
$ws.Range("I5").Value = 0.1531387528565491
$ws.Range("J5").Value = 0.1531387528565491
$ws.Range("M5").Value = 14.82137566666667
$ws.Range("N5").Value = 44.464127
$ws.Range("O5").Value = 0.9193095050964014
$ws.Range("P5").Value = 0.9193095050964012
$ws.Range("Q5").Value = 310.3903212908399
$ws.Range("R5").Value = 2793.512891617559
$ws.Range("S5").Value = 0.1407819110996343
$ws.Range("T5").Value = 0.1407819110996343

# Row 6
$ws.Range("I6").Value = 0.1531387528565491
$ws.Range("J6").Value = 0.1531387528565491
$ws.Range("O6").Value = 0.0119476607067887
$ws.Range("P6").Value = 0.0119476607067887
$ws.Range("S6").Value = 0.001829649860190818
$ws.Range("T6").Value = 0.001829649860190817

# Row 7
$ws.Range("I7").Value = 0.1531387528565491
$ws.Range("J7").Value = 0.1531387528565491
$ws.Range("M7").Value = 1.108292
$ws.Range("N7").Value = 3.324876
$ws.Range("O7").Value = 0.06874283419681
$ws.Range("P7").Value = 0.06874283419680997
$ws.Range("Q7").Value = 23.20993123045467
$ws.Range("R7").Value = 208.889381074092
$ws.Range("S7").Value = 0.01052719189672402
$ws.Range("T7").Value = 0.01052719189672401

# Row 8
$ws.Range("G8").Value = 2.220865666666667
$ws.Range("H8").Value = 6.662597
$ws.Range("I8").Value = 0.01624006416566169
$ws.Range("J8").Value = 0.01624006416566169
$ws.Range("M8").Value = 14.82137566666667
$ws.Range("N8").Value = 44.464127
$ws.Range("O8").Value = 0.9193095050964014
$ws.Range("P8").Value = 0.9193095050964012
$ws.Range("Q8").Value = 32.91628435086878
$ws.Range("R8").Value = 296.246559157819
$ws.Range("S8").Value = 0.01492964535086826
$ws.Range("T8").Value = 0.01492964535086825

# Row 9
$ws.Range("G9").Value = 2.220865666666667
$ws.Range("H9").Value = 6.662597
$ws.Range("I9").Value = 0.01624006416566169
$ws.Range("J9").Value = 0.01624006416566169
$ws.Range("O9").Value = 0.0119476607067887
$ws.Range("P9").Value = 0.0119476607067887
$ws.Range("Q9").Value = 0.4277912878874444
$ws.Range("R9").Value = 3.850121590987
$ws.Range("S9").Value = 0.0001940307765078034
$ws.Range("T9").Value = 0.0001940307765078034

# Row 10
$ws.Range("G10").Value = 2.220865666666667
$ws.Range("H10").Value = 6.662597
$ws.Range("I10").Value = 0.01624006416566169
$ws.Range("J10").Value = 0.01624006416566169
$ws.Range("M10").Value = 1.108292
$ws.Range("N10").Value = 3.324876
$ws.Range("O10").Value = 0.06874283419681
$ws.Range("P10").Value = 0.06874283419680997
$ws.Range("Q10").Value = 2.461367651441333
$ws.Range("R10").Value = 22.152308862972
$ws.Range("S10").Value = 0.001116388038285638
$ws.Range("T10").Value = 0.001116388038285637
